$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Strike-through three existing "done" to-do items (paragraphs
#    1, 3 and 4 - the hyperlink item, the "Add tabs" item and the
#    "Change the script for the time-picker" item). Setting
#    Font.StrikeThrough on the paragraph's Range (which includes the
#    paragraph mark) stamps <w:strike/> onto every run's rPr as well
#    as onto the paragraph mark's rPr, matching the diff exactly.
# ------------------------------------------------------------------
$d.Paragraphs.Item(1).Range.Font.StrikeThrough = 1
$d.Paragraphs.Item(3).Range.Font.StrikeThrough = 1
$d.Paragraphs.Item(4).Range.Font.StrikeThrough = 1

# ------------------------------------------------------------------
# 2) Last paragraph ("Enable up" / bookmark / "loader to be able...")
#    gets merged into a single run with the full sentence, and the
#    _GoBack bookmark is removed from there (it is re-created at the
#    end of the new "Adjust menubar alignment" paragraph below).
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIdx)
$fullRange = $d.Range($pLast.Range.Start, $pLast.Range.End)
$fullRange.Text = "Enable uploader to be able to circle a person in the photo or able to fuzz out people not related."

# ------------------------------------------------------------------
# 3) Insert two new to-do items after it:
#      "Make full page background white"                  (struck out)
#      "Adjust menubar alignment" (with spell-check marks) (struck out)
#    and put the _GoBack bookmark back at the very end of the
#    document (end of the last paragraph).
# ------------------------------------------------------------------
$pLast = $d.Paragraphs.Item($d.Paragraphs.Count)
$pLast.Range.InsertParagraphAfter()

$pNew1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pNew1.Range.Text = "Make full page background white"
$pNew1.Range.Font.StrikeThrough = 1

$pNew1 = $d.Paragraphs.Item($d.Paragraphs.Count)
$pNew1.Range.InsertParagraphAfter()

$pNew2 = $d.Paragraphs.Item($d.Paragraphs.Count)

$xmlFrag = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Adjust </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>menubar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> alignment</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$pNew2.Range.InsertXML($xmlFrag)

$pNew2 = $d.Paragraphs.Item($d.Paragraphs.Count)
$d.Bookmarks.Add("_GoBack", $pNew2.Range.Characters.Last) | Out-Null

Write-Host "Done. Paragraph count:" $d.Paragraphs.Count
